# Random forest and ensembling
# Applies the ModelTracking.xlsx update: drops the now-redundant "None"
# hyperparameter placeholders from the older rows, tidies up column B's
# formatting, and appends the Random Forest + ensembling experiment rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# ---------------------------------------------------------------------
# 1. Drop the placeholder "None" hyperparameter values from rows 2-6 -
#    those models were run with default params, so the column is blank now.
# ---------------------------------------------------------------------
$ws.Range("C2:C6").ClearContents()

# ---------------------------------------------------------------------
# 2. Column B (Date) - strip the stray "applyNumberFormat" flag that was
#    sitting on the column default, keep the header bold/General and the
#    data cells formatted as dates.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ClearFormats()
$ws.Range("B1").Font.Bold = $true
$ws.Range("B2:B8").NumberFormat = "mmm-yy"

# ---------------------------------------------------------------------
# 3. Column C gets wider now that it holds hyperparameter notes; column H
#    (Details) gets a lot wider too - both wrap their text where used.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 27.8
$ws.Columns.Item(8).ColumnWidth = 44.95

# Row 7 (Decision Tree Regressor) - existing data, just needs the Details
# cell to wrap and the taller row that comes with it.
$ws.Range("H7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 28.8

# ---------------------------------------------------------------------
# 4. New row 8 - Random Forest Regressor results.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Random Forest Regessor"
$ws.Range("B8").Value = 46447
$ws.Range("C8").Value = "Random state = 61" + $nl + "n_jobs = -1"
$ws.Range("C8").WrapText = $true
$ws.Range("D8").Value = 474.37886939340001
$ws.Range("E8").Value = 1367.15354958689
$ws.Range("H8").Value = "This Prouces the best result as of now and can be tuned for better accuracy"
$ws.Range("H8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 28.8

# ---------------------------------------------------------------------
# 5. New row 9 - ensembled (80% Random Forest / 20% Ridge) results.
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "80 - Random Forest" + $nl + "20 - Ridge"
$ws.Range("A9").WrapText = $true
$ws.Range("B9").Value = 46447
$ws.Range("B9").NumberFormat = "d-mmm"
$ws.Range("C9").Value = "n_jobs=-1, random_state=len(x_train), n_estimators=20" + $nl + "-----" + $nl + "None"
$ws.Range("C9").WrapText = $true
$ws.Range("D9").Value = 812.04496192446902
$ws.Range("E9").Value = 1431.3805071645099
$ws.Range("H9").Value = "The combination is better than many of the models but not more than that of the Random forest's Original Output"
$ws.Range("H9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 72

# Header "Details" also wraps now that column H is so much wider.
$ws.Range("H1").WrapText = $true

# ---------------------------------------------------------------------
# 6. Selection follows the last entry row, and print orientation is set
#    to portrait (as recorded by the page setup dialog).
# ---------------------------------------------------------------------
$ws.Range("A10").Select()
$ws.PageSetup.Orientation = 1
